$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Reverse the "Periodo Mora" values in E16:E22 (2308..2402 -> 2402..2308)
$periodos = @("2402", "2401", "2312", "2311", "2310", "2309", "2308")
for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periodos[$i]
}

# Swap the "Valor Mora" amounts between F16 and F22
$ws.Range("F16").Value = 68400
$ws.Range("F22").Value = 114000
